$d = $word.ActiveDocument
$count = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($count)
$last.Range.InsertParagraphAfter()
$count2 = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($count2)
$xml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Requisitos no funcionales</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t>Usabilidad</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Interfaz intuitiva con Bootstrap 5</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Tiempo de respuesta menor a 3 segundos</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Navegación consistente entre módulos</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t>Rendimiento</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Soporte para múltiples usuarios concurrentes</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Capacidad de manejar 1000+ ventas diarias</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Backup automático diario de base de datos</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t>Seguridad</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Autenticación de usuarios</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Control de acceso por roles</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Registro de auditoría de operaciones</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t>Confiabilidad</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Disponibilidad 24/7</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Tolerancia a fallos en operaciones críticas</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Validación de integridad de datos</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $newPara.Range.InsertXML($xml)

# InsertXML leaves behind the placeholder paragraph's own (now empty) mark
# as a trailing empty paragraph. Merge it away by deleting the range that
# spans from just before it (the end of the real last paragraph) through
# its own end.
$finalCount = $d.Paragraphs.Count
$delStart = $d.Paragraphs.Item($finalCount - 1).Range.End - 1
$delEnd = $d.Paragraphs.Item($finalCount).Range.End
$null = $d.Range($delStart, $delEnd).Delete()
